$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 445, shifting existing rows 445:544 down to 446:545.
$ws.Rows.Item(445).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A445").Value = 6
$ws.Range("B445").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C445").Value = "Metropolitana"
$ws.Range("D445").Value = 44889
$ws.Range("E445").Value = 13
$ws.Range("F445").Value = 100112052
$ws.Range("G445").Value = "Albahaca"
$ws.Range("H445").Value = "Sin especificar"
$ws.Range("I445").Value = "Primera"
$ws.Range("J445").Value = 530
$ws.Range("K445").Value = 5000
$ws.Range("L445").Value = 5500
$ws.Range("M445").Value = 5283
$ws.Range("N445").Value = "`$/docena de matas"
$ws.Range("O445").Value = "Región Metropolitana"
$ws.Range("P445").Value = 880
$ws.Range("Q445").Value = 6
$ws.Range("R445").Value = "Hortaliza"
